$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: now computed as the midpoint between A2 and A3 (as fractions)
$ws.Range("D2").Formula = "=(A2/100+(A3/100-A2/100)/2)"

# F2: now the running total of the Q column
$ws.Range("F2").Formula = "=SUM(E2:E12)"

# D3 / E3: midpoint + incremental-Q formulas
$ws.Range("D3").Formula = "=(A3/100+(A4/100-A3/100)/2)"
$ws.Range("E3").Formula = "=(D3-D2)*(B3/100)*C3"

# D4:D8 / E4:E8: same formulas filled down (shared formula ranges)
$ws.Range("D4:D8").Formula = "=(A4/100+(A5/100-A4/100)/2)"
$ws.Range("E4:E8").Formula = "=(D4-D3)*(B4/100)*C4"

# These cells lose their explicit style (s="1") and revert to the default style
$ws.Range("D2:D8").Style = "Normal"
$ws.Range("E3:E8").Style = "Normal"
$ws.Range("F2").Style = "Normal"

# Update the active selection
$ws.Range("E15").Select()
